# Update with Correct Forecast output
# - Insert a "Week_Start_Date" column (new col B) on "Forecast Comparison"
# - Shorten week labels in col A (W01 -> W1, ... W16 stays W16)
# - Fix two MyForecast values (W3: 42 -> 41, W16: 37 -> 36)
# - Store is_holiday_week as boolean (FALSE) instead of numeric 0
# - Update dependent totals/min on the "Summary" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# 1) Insert new column B ("Week_Start_Date"), shifting ASIN/MyForecast/etc. right.
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "Week_Start_Date"

# Week start dates (Mondays) for the 16 forecast weeks. Keep them as plain
# text (not real dates) to match the source data's string representation.
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    # Shorten "W01".."W16" -> "W1".."W16"
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
    # is_holiday_week now stored as a boolean rather than 0/1 numeric.
    $ws.Cells.Item($row, 10).Value = $false
}

# 2) Corrected MyForecast values (column D after the insert).
$ws.Cells.Item(4, 4).Value = 41
$ws.Cells.Item(17, 4).Value = 36

# 3) Propagate the corrected totals/min to the Summary sheet (kept as text,
#    matching the sheet's existing string-valued "Value" column).
$summary.Range("B9:B11").NumberFormat = "@"
$summary.Range("B14").NumberFormat = "@"
$summary.Range("B9").Value = "699"
$summary.Range("B10").Value = "340"
$summary.Range("B11").Value = "169"
$summary.Range("B14").Value = "36"
